{"js": "// Replace each two-digit multiplication expression with its new value.\n// Each \"before\" text is unique in the document, so a simple search +\n// replace (one hit per pair) is unambiguous and safe.\nconst replacements = [\n  [\"69\u00d726=\", \"94\u00d724=\"],\n  [\"89\u00d733=\", \"98\u00d792=\"],\n  [\"47\u00d727=\", \"91\u00d786=\"],\n  [\"38\u00d735=\", \"13\u00d763=\"],\n  [\"60\u00d769=\", \"30\u00d755=\"],\n  [\"77\u00d719=\", \"38\u00d746=\"],\n  [\"33\u00d787=\", \"87\u00d728=\"],\n  [\"72\u00d728=\", \"86\u00d759=\"],\n  [\"40\u00d771=\", \"51\u00d733=\"],\n  [\"64\u00d799=\", \"99\u00d791=\"],\n  [\"17\u00d775=\", \"52\u00d773=\"],\n  [\"78\u00d744=\", \"27\u00d745=\"],\n  [\"64\u00d781=\", \"91\u00d791=\"],\n  [\"93\u00d727=\", \"49\u00d741=\"],\n  [\"74\u00d780=\", \"58\u00d764=\"],\n  [\"84\u00d724=\", \"96\u00d799=\"],\n  [\"87\u00d722=\", \"51\u00d777=\"],\n  [\"52\u00d771=\", \"43\u00d772=\"],\n  [\"34\u00d729=\", \"44\u00d786=\"],\n  [\"85\u00d790=\", \"38\u00d779=\"],\n  [\"94\u00d748=\", \"23\u00d760=\"],\n  [\"72\u00d714=\", \"82\u00d774=\"],\n  [\"70\u00d751=\", \"39\u00d760=\"],\n  [\"65\u00d781=\", \"76\u00d723=\"],\n  [\"38\u00d751=\", \"85\u00d770=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each two-digit multiplication expression with its new value.\n# Each \"before\" text is unique in the document, so Find/Replace with\n# wdReplaceAll (which will only ever touch the single matching occurrence)\n# is unambiguous and safe for every pair.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"69\u00d726=\", \"94\u00d724=\"),\n    @(\"89\u00d733=\", \"98\u00d792=\"),\n    @(\"47\u00d727=\", \"91\u00d786=\"),\n    @(\"38\u00d735=\", \"13\u00d763=\"),\n    @(\"60\u00d769=\", \"30\u00d755=\"),\n    @(\"77\u00d719=\", \"38\u00d746=\"),\n    @(\"33\u00d787=\", \"87\u00d728=\"),\n    @(\"72\u00d728=\", \"86\u00d759=\"),\n    @(\"40\u00d771=\", \"51\u00d733=\"),\n    @(\"64\u00d799=\", \"99\u00d791=\"),\n    @(\"17\u00d775=\", \"52\u00d773=\"),\n    @(\"78\u00d744=\", \"27\u00d745=\"),\n    @(\"64\u00d781=\", \"91\u00d791=\"),\n    @(\"93\u00d727=\", \"49\u00d741=\"),\n    @(\"74\u00d780=\", \"58\u00d764=\"),\n    @(\"84\u00d724=\", \"96\u00d799=\"),\n    @(\"87\u00d722=\", \"51\u00d777=\"),\n    @(\"52\u00d771=\", \"43\u00d772=\"),\n    @(\"34\u00d729=\", \"44\u00d786=\"),\n    @(\"85\u00d790=\", \"38\u00d779=\"),\n    @(\"94\u00d748=\", \"23\u00d760=\"),\n    @(\"72\u00d714=\", \"82\u00d774=\"),\n    @(\"70\u00d751=\", \"39\u00d760=\"),\n    @(\"65\u00d781=\", \"76\u00d723=\"),\n    @(\"38\u00d751=\", \"85\u00d770=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n\n    $find.Execute([ref]$oldText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$newText, [ref]2)\n}\n"}
